$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.579.01"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "2.287.36"
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'313.82"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").Value = "'105.08"
$ws.Range("E6").Value = "  +2.41%  "
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "'0.605"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").Value = "'39.59"
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("D11").Value = "'0.0905"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "'8.39"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("E13").Value = "  +2.68%  "
$ws.Range("E14").Value = "  +3.58%  "
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "2.636.01"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "2.283.20"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "42.775.36"
$ws.Range("D19").Value = "'7.42"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.0000106"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").Value = "'13.52"
$ws.Range("E21").Value = "  +19.18%  "
$ws.Range("D22").Value = "'73.98"
$ws.Range("E23").Value = "  +1.16%  "
$ws.Range("D24").Value = "'265.18"
$ws.Range("E24").Value = "  -4.22%  "
$ws.Range("E25").Value = "  -1.79%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").Value = "'10.83"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "'7.21"
$ws.Range("E28").Value = "  +23.09%  "
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("D31").Value = "'37.29"
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("D33").Value = "'0.0877"
$ws.Range("E34").Value = "  -2.89%  "
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("D36").Value = "'0.114"
$ws.Range("E36").Value = "  -3.72%  "
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("E38").Value = "  -3.79%  "
$ws.Range("D39").Value = "'3.78"
$ws.Range("E39").Value = "  +2.54%  "
$ws.Range("D40").Value = "'2.66"
$ws.Range("E40").Value = "  -3.64%  "
$ws.Range("D41").Value = "'1.57"
$ws.Range("E41").Value = "  +4.82%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.233"
$ws.Range("E42").Value = "  +3.27%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").Value = "'70.71"
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D45").Value = "'94.34"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("D46").Value = "'12.14"
$ws.Range("E46").Value = "  +1.21%  "
$ws.Range("D47").Value = "1.738.27"
$ws.Range("E47").Value = "  +9.42%  "
$ws.Range("D48").Value = "'112.86"
$ws.Range("E48").Value = "  +0.18%  "
$ws.Range("D49").Value = "'79.96"
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'8.76"
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'5.21"
$ws.Range("E51").Value = "  -0.52%  "
